# Week-four description largely done:
#  - Rename the active test-set sheet from "4.3.3 (Feb 2nd)" to "4.3.3 (Jan 27th)"
#  - Move/resize the saved window view
#  - Update the remembered selection on the active sheet

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "4.3.3 (Jan 27th)"

# Reposition / resize the application window the way it was left after editing.
$win = $excel.Windows.Item(1)
$win.Left = 4600
$win.Top = 22100
$win.Width = 28800
$win.Height = 17500

# Re-select the cell the author ended up on.
$ws1.Activate()
$ws1.Range("C32").Select()
